$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(1,2,3,4,5,6,7,8,9)
for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 2  # B=2 ... J=10
    $ws.Cells.Item(1, $col).Value = $values[$i]
}
